$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 5
$ws.Range("H5").Value = 352.29413
$ws.Range("I5").Value = 221
$ws.Range("K5").Value = 221
$ws.Range("M5").Value = -106
# Row 6
$ws.Range("H6").Value = 3170.7058
$ws.Range("I6").Value = 134.66667
$ws.Range("J6").Value = 3821.2856
$ws.Range("K6").Value = 404.00001
$ws.Range("L6").Value = 11463.8568
$ws.Range("M6").Value = -292.00001
$ws.Range("N6").Value = -11687.8568
# Row 12
$ws.Range("H12").Value = 1497
$ws.Range("I12").Value = 1796.1428
$ws.Range("J12").Value = 450
$ws.Range("K12").Value = 1796.1428
$ws.Range("L12").Value = 450
$ws.Range("M12").Value = -1626.1428
$ws.Range("N12").Value = -790
# Row 29
$ws.Range("H29").Value = 1077.2941
$ws.Range("I29").Value = 862.75
$ws.Range("J29").Value = 1268
$ws.Range("K29").Value = 2588.25
$ws.Range("L29").Value = 3804
$ws.Range("M29").Value = -2307.25
$ws.Range("N29").Value = -4366
# Row 38
$ws.Range("H38").Value = 3263.3
$ws.Range("I38").Value = 199.42857
$ws.Range("J38").Value = 4913.077
$ws.Range("K38").Value = 598.28571
$ws.Range("L38").Value = 14739.231
$ws.Range("M38").Value = -226.28571
$ws.Range("N38").Value = -15483.231
# Row 58
$ws.Range("H58").Value = 3944.889
$ws.Range("I58").Value = 688
$ws.Range("J58").Value = 30000
$ws.Range("K58").Value = 2064
$ws.Range("L58").Value = 90000
$ws.Range("M58").Value = -1914
$ws.Range("N58").Value = -90300
# Row 64
$ws.Range("H64").Value = 3492.4614
$ws.Range("I64").Value = 3114.5715
$ws.Range("J64").Value = 3933.3333
$ws.Range("K64").Value = 3114.5715
$ws.Range("L64").Value = 3933.3333
$ws.Range("M64").Value = -2866.5715
$ws.Range("N64").Value = -4429.3333
# Row 67
$ws.Range("H67").Value = 3492.4614
$ws.Range("I67").Value = 3114.5715
$ws.Range("J67").Value = 3933.3333
$ws.Range("K67").Value = 3114.5715
$ws.Range("L67").Value = 3933.3333
$ws.Range("M67").Value = -2256.5715
$ws.Range("N67").Value = -5649.3333
# Row 129
$ws.Range("H129").Value = 845.92
$ws.Range("J129").Value = 866.61456
$ws.Range("L129").Value = 2599.84368
$ws.Range("N129").Value = -12599.84368
# Row 132
$ws.Range("H132").Value = 25004158
$ws.Range("I132").Value = 26319692
$ws.Range("J132").Value = 9000
$ws.Range("K132").Value = 78959076
$ws.Range("L132").Value = 27000
$ws.Range("M132").Value = -78956546
$ws.Range("N132").Value = -32060

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 4
$ws.Range("H4").Value = 200
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 200
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 200
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = -432
# Row 23
$ws.Range("H23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()
# Row 45
$ws.Range("H45").Value = 3745
$ws.Range("I45").Value = 4326.6665
$ws.Range("K45").Value = 4326.6665
$ws.Range("M45").Value = -3949.6665
# Row 61
$ws.Range("H61").Value = 1587.091
$ws.Range("I61").Value = 1606.4445
$ws.Range("J61").Value = 1500
$ws.Range("K61").Value = 1606.4445
$ws.Range("L61").Value = 1500
$ws.Range("M61").Value = -1394.4445
$ws.Range("N61").Value = -1924
# Row 74
$ws.Range("H74").Value = 4378.7856
$ws.Range("I74").Value = 4289.222
$ws.Range("K74").Value = 4289.222
$ws.Range("M74").Value = -3415.222
# Row 77
$ws.Range("H77").Value = 4378.7856
$ws.Range("I77").Value = 4289.222
$ws.Range("K77").Value = 21446.11
$ws.Range("M77").Value = -17078.11
# Row 125
$ws.Range("H125").Value = 41669
$ws.Range("J125").Value = 41669
$ws.Range("L125").Value = 41669
$ws.Range("N125").Value = -51509
# Row 132
$ws.Range("H132").Value = 2655.45
$ws.Range("I132").Value = 1258.9166
$ws.Range("J132").Value = 4750.25
$ws.Range("K132").Value = 3776.7498
$ws.Range("L132").Value = 14250.75
$ws.Range("M132").Value = -1246.7498
$ws.Range("N132").Value = -19310.75
# Row 136
$ws.Range("H136").Value = 1587.091
$ws.Range("I136").Value = 1606.4445
$ws.Range("J136").Value = 1500
$ws.Range("K136").Value = 4819.333500000001
$ws.Range("L136").Value = 4500
$ws.Range("M136").Value = -2269.333500000001
$ws.Range("N136").Value = -9600
# Row 139
$ws.Range("H139").Value = 41663.395
$ws.Range("J139").Value = 41663.395
$ws.Range("L139").Value = 41663.395
$ws.Range("N139").Value = -51943.395

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 22
$ws.Range("H22").Value = 378.42856
$ws.Range("I22").Value = 357.83334
$ws.Range("J22").Value = 502
$ws.Range("K22").Value = 357.83334
$ws.Range("L22").Value = 502
$ws.Range("M22").Value = -184.83334
$ws.Range("N22").Value = -848

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 742.8570999999999
$ws.Range("I22").Value = 438.83334
$ws.Range("J22").Value = 1290.1
$ws.Range("K22").Value = 438.83334
$ws.Range("L22").Value = 1290.1
$ws.Range("M22").Value = -88.83334000000002
$ws.Range("N22").Value = -1990.1
# Row 23
$ws.Range("H23").Value = 25750
$ws.Range("I23").Value = 16500
$ws.Range("K23").Value = 16500
$ws.Range("M23").Value = -16260
# Row 27
$ws.Range("H27").Value = 25750
$ws.Range("I27").Value = 16500
$ws.Range("K27").Value = 16500
$ws.Range("M27").Value = -16308
# Row 31
$ws.Range("H31").Value = 5884.0513
$ws.Range("I31").Value = 2447.7693
$ws.Range("J31").Value = 12756.615
$ws.Range("K31").Value = 2447.7693
$ws.Range("L31").Value = 12756.615
$ws.Range("M31").Value = -2152.7693
$ws.Range("N31").Value = -13346.615
# Row 34
$ws.Range("H34").Value = 5884.0513
$ws.Range("I34").Value = 2447.7693
$ws.Range("J34").Value = 12756.615
$ws.Range("K34").Value = 2447.7693
$ws.Range("L34").Value = 12756.615
$ws.Range("M34").Value = -2245.7693
$ws.Range("N34").Value = -13160.615
# Row 52
$ws.Range("H52").Value = 34600
$ws.Range("J52").Value = 34600
$ws.Range("L52").Value = 34600
$ws.Range("N52").Value = -35188
# Row 99
$ws.Range("H99").Value = 11115393
$ws.Range("I99").Value = 33335366
$ws.Range("J99").Value = 5407.5
$ws.Range("K99").Value = 33335366
$ws.Range("L99").Value = 5407.5
$ws.Range("M99").Value = -33333868
$ws.Range("N99").Value = -8403.5
# Row 126
$ws.Range("H126").Value = 11115393
$ws.Range("I126").Value = 33335366
$ws.Range("J126").Value = 5407.5
$ws.Range("K126").Value = 100006098
$ws.Range("L126").Value = 16222.5
$ws.Range("M126").Value = -100003628
$ws.Range("N126").Value = -21162.5
# Row 132
$ws.Range("H132").Value = 4790.6665
$ws.Range("I132").Value = 2026.8
$ws.Range("J132").Value = 8245.5
$ws.Range("K132").Value = 6080.4
$ws.Range("L132").Value = 24736.5
$ws.Range("M132").Value = -3550.4
$ws.Range("N132").Value = -29796.5
# Row 137
$ws.Range("H137").Value = 44557.145
$ws.Range("J137").Value = 44557.145
$ws.Range("L137").Value = 44557.145
$ws.Range("N137").Value = -54757.145

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 3
$ws.Range("H3").Value = 2892.2222
$ws.Range("I3").Value = 1938.3334
$ws.Range("K3").Value = 5815.0002
$ws.Range("M3").Value = -5703.0002
# Row 12
$ws.Range("H12").Value = 89.21738999999999
$ws.Range("I12").Value = 23.142857
$ws.Range("J12").Value = 118.125
$ws.Range("K12").Value = 69.42857100000001
$ws.Range("L12").Value = 354.375
$ws.Range("M12").Value = 103.571429
$ws.Range("N12").Value = -700.375
# Row 129
$ws.Range("H129").Value = 2760.353
$ws.Range("J129").Value = 2295.0908
$ws.Range("L129").Value = 6885.2724
$ws.Range("N129").Value = -16885.2724
# Row 131
$ws.Range("H131").Value = 829.14
$ws.Range("J131").Value = 850.6667
$ws.Range("L131").Value = 2552.0001
$ws.Range("N131").Value = -12632.0001

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 99
$ws.Range("H99").Value = 20000
$ws.Range("I99").Value = 20000
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 20000
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -17754
$ws.Range("N99").ClearContents()
# Row 132
$ws.Range("H132").Value = 3139.5417
$ws.Range("I132").Value = 2292.55
$ws.Range("K132").Value = 6877.650000000001
$ws.Range("M132").Value = -4347.650000000001
# Row 140
$ws.Range("H140").Value = 38728.65
$ws.Range("J140").Value = 38728.65
$ws.Range("L140").Value = 38728.65
$ws.Range("N140").Value = -49088.65

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 5
$ws.Range("H5").Value = 30200
$ws.Range("I5").Value = 30400
$ws.Range("J5").Value = 30000
$ws.Range("K5").Value = 30400
$ws.Range("L5").Value = 30000
$ws.Range("M5").Value = -30287
$ws.Range("N5").Value = -30226
# Row 40
$ws.Range("H40").Value = 9100.333000000001
$ws.Range("I40").Value = 7557.5713
$ws.Range("J40").Value = 14500
$ws.Range("K40").Value = 7557.5713
$ws.Range("L40").Value = 14500
$ws.Range("M40").Value = -7421.5713
$ws.Range("N40").Value = -14772
# Row 93
$ws.Range("H93").Value = 2661.2666
$ws.Range("J93").Value = 3520.4
$ws.Range("L93").Value = 3520.4
$ws.Range("N93").Value = -6016.4

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 21
$ws.Range("H21").Value = 40000
$ws.Range("J21").Value = 40000
$ws.Range("L21").Value = 40000
$ws.Range("N21").Value = -40470
# Row 22
$ws.Range("H22").Value = 23674.75
$ws.Range("J22").Value = 23674.75
$ws.Range("L22").Value = 23674.75
$ws.Range("N22").Value = -24260.75
# Row 35
$ws.Range("H35").Value = 40000
$ws.Range("J35").Value = 40000
$ws.Range("L35").Value = 40000
$ws.Range("N35").Value = -40580
# Row 46
$ws.Range("H46").Value = 45393
$ws.Range("J46").Value = 45393
$ws.Range("L46").Value = 45393
$ws.Range("N46").Value = -45855
# Row 81
$ws.Range("H81").Value = 1833.3334
$ws.Range("I81").Value = 1833.3334
$ws.Range("K81").Value = 3666.6668
$ws.Range("M81").Value = -2605.6668
# Row 84
$ws.Range("H84").Value = 1833.3334
$ws.Range("I84").Value = 1833.3334
$ws.Range("K84").Value = 18333.334
$ws.Range("M84").Value = -13029.334
# Row 134
$ws.Range("H134").Value = 45393
$ws.Range("J134").Value = 45393
$ws.Range("L134").Value = 136179
$ws.Range("N134").Value = -141249
